# NIT-8000699336.xlsx — "Elimina EC anteriores y se agregan nuevos, se modifica base de datos"
#
# Updates the account-statement (Estado de Cuenta) worksheet:
#  - bumps the two summary numbers (Valor Mora total, Cant. Trabajadores, Cant. Periodos)
#  - replaces the worker/period detail table (rows 16-28) with a new table
#    (rows 16-31) containing a new worker (KATERIN YULIETH LAGOS BELLO) plus
#    refreshed period rows (2501-2507/2412) for the two existing workers, all
#    at the new Salario Basico of 1,300,000 / 1,423,500.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Make room: the detail table grows from 13 rows (16-28) to 16 rows
#    (16-31). Insert 3 blank rows right after the old last data row (28) —
#    this also pushes the footer block (old rows 33/34) down to 36/37,
#    exactly matching the target layout.
# ---------------------------------------------------------------------------
$ws.Range("A29:A31").EntireRow.Insert()

# ---------------------------------------------------------------------------
# 2. Fix up formatting of the (old + new) last few data rows.
#    Row 28 currently still carries the special "bottom of table" style
#    (thicker border); that style needs to move down to the new last row
#    (31), while rows 28-30 become normal interior rows.
# ---------------------------------------------------------------------------
# Capture the "bottom of table" formatting (currently on row 28) onto the new
# last row (31) before row 28's own formatting gets overwritten below.
$ws.Range("B28:J28").Copy()
$ws.Range("B31:J31").PasteSpecial(-4122)

# Re-stripe rows 28-30 with the normal interior-row formatting (taken from
# row 27, an untouched normal row).
$ws.Range("B27:J27").Copy()
$ws.Range("B28:J28").PasteSpecial(-4122)
$ws.Range("B29:J29").PasteSpecial(-4122)
$ws.Range("B30:J30").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 3. Rewrite the detail table contents, rows 16-31.
# ---------------------------------------------------------------------------
$tipoDoc = "CC"

$data = @(
  @("1193554162", "KATERIN YULIETH LAGOS BELLO",     "2507", 56940, 1423500),
  @("1143381131", "BELKYS CAROLINA RICARDO ROMERO",  "2507", 52000, 1300000),
  @("1143381131", "BELKYS CAROLINA RICARDO ROMERO",  "2506", 52000, 1300000),
  @("1143381131", "BELKYS CAROLINA RICARDO ROMERO",  "2505", 52000, 1300000),
  @("1143381131", "BELKYS CAROLINA RICARDO ROMERO",  "2504", 52000, 1300000),
  @("1143381131", "BELKYS CAROLINA RICARDO ROMERO",  "2503", 52000, 1300000),
  @("1143381131", "BELKYS CAROLINA RICARDO ROMERO",  "2502", 52000, 1300000),
  @("1143381131", "BELKYS CAROLINA RICARDO ROMERO",  "2501", 52000, 1300000),
  @("1143381131", "BELKYS CAROLINA RICARDO ROMERO",  "2412", 43333, 1300000),
  @("1143349310", "BENISPAULET VILLAR ACEVEDO",      "2507", 52000, 1300000),
  @("1143349310", "BENISPAULET VILLAR ACEVEDO",      "2506", 52000, 1300000),
  @("1143349310", "BENISPAULET VILLAR ACEVEDO",      "2505", 52000, 1300000),
  @("1143349310", "BENISPAULET VILLAR ACEVEDO",      "2504", 52000, 1300000),
  @("1143349310", "BENISPAULET VILLAR ACEVEDO",      "2503", 52000, 1300000),
  @("1143349310", "BENISPAULET VILLAR ACEVEDO",      "2502", 52000, 1300000),
  @("1143349310", "BENISPAULET VILLAR ACEVEDO",      "2501", 52000, 1300000)
)

$r = 16
foreach ($row in $data) {
    $ws.Cells.Item($r, 2).Value = $tipoDoc
    $ws.Cells.Item($r, 3).Value = $row[0]
    $ws.Cells.Item($r, 4).Value = $row[1]
    $ws.Cells.Item($r, 5).Value = $row[2]
    $ws.Cells.Item($r, 6).Value = $row[3]
    $ws.Cells.Item($r, 7).Value = $row[4]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# 4. Update the summary fields above the table.
# ---------------------------------------------------------------------------
$ws.Range("E11").Value = 828273     # VALOR MORA (total)
$ws.Range("C13").Value = 3          # Cant. Trabajadores
$ws.Range("F13").Value = 8          # Cant. Periodos
